$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "L120.csv"
$ws.Range("B6").Value = "RMS"
$ws.Range("C6").Value = "grid1"
$ws.Range("D6").Value = 4.9837033567055027
$ws.Range("E6").Value = 7.6154278127341932
$ws.Range("F6").Value = -0.42355894430186491
$ws.Range("G6").Value = -0.39515686533869077
$ws.Range("H6").Value = 0.37708656173945287
$ws.Range("I6").Value = 2.444139291155317
$ws.Range("J6").Value = 4.4169010412501262
$ws.Range("K6").Value = "grid2"
$ws.Range("L6").Value = 4.8132933006222851
$ws.Range("M6").Value = 26.497697273869736
$ws.Range("N6").Value = -0.25380222280809078
$ws.Range("O6").Value = -0.40754931110034676
$ws.Range("P6").Value = 0.55743954844668608
$ws.Range("Q6").Value = 2.6769939795094224
$ws.Range("R6").Value = 4.8598205241911741
